$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.118.40'
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").Value = '1.787.02'
$ws.Range("E3").Value = '  -0.18%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").Value = '''226.83'
$ws.Range("E5").Value = '  -0.87%  '

$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("D8").Value = '''31.94'
$ws.Range("E8").Value = '  -1.74%  '

$ws.Range("E9").Value = '  +1.28%  '

$ws.Range("E10").Value = '  -3.03%  '

$ws.Range("D11").Value = '''0.0945'
$ws.Range("E11").Value = '  +0.93%  '

$ws.Range("E12").Value = '  -0.20%  '

$ws.Range("E13").Value = '  +1.35%  '

$ws.Range("D14").Value = '1.768.06'
$ws.Range("E14").Value = '  -1.27%  '

$ws.Range("D15").Value = '34.030.36'
$ws.Range("E15").Value = '  -0.11%  '

$ws.Range("D16").Value = '''0.620'
$ws.Range("E16").Value = '  -0.90%  '

$ws.Range("D17").Value = '''4.18'
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("D18").Value = '''68.01'
$ws.Range("E18").Value = '  -0.69%  '

$ws.Range("D19").Value = '''245.24'
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("E20").Value = '  -1.53%  '

$ws.Range("E21").Value = '  +0.25%  '

$ws.Range("D22").Value = '''10.85'
$ws.Range("E22").Value = '  +0.82%  '

$ws.Range("E23").Value = '  -0.36%  '

$ws.Range("E24").Value = '  -2.29%  '

$ws.Range("D25").Value = '''161.41'
$ws.Range("E25").Value = '  +0.41%  '

$ws.Range("E26").Value = '  +0.50%  '

$ws.Range("D27").Value = '''16.30'
$ws.Range("E27").Value = '  -0.49%  '

$ws.Range("E28").Value = '  +0.20%  '

$ws.Range("E29").Value = '  +0.30%  '

$ws.Range("E30").Value = '  -1.59%  '

$ws.Range("E31").Value = '  +0.14%  '

$ws.Range("D32").Value = '''3.66'
$ws.Range("E32").Value = '  -0.74%  '

$ws.Range("D33").Value = '''3.61'
$ws.Range("E33").Value = '  +2.28%  '

$ws.Range("D34").Value = '''1.81'
$ws.Range("E34").Value = '  -0.91%  '

$ws.Range("D35").Value = '1.454.12'
$ws.Range("E35").Value = '  +3.95%  '

$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").Value = '''2.42'
$ws.Range("E36").Value = '  +9.34%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.646'
$ws.Range("E37").Value = '  -2.40%  '

$ws.Range("E38").Value = '  +1.79%  '

$ws.Range("D39").Value = '''1.04'
$ws.Range("E39").Value = '  -1.05%  '

$ws.Range("D40").Value = '''80.20'
$ws.Range("E40").Value = '  +2.24%  '

$ws.Range("E41").Value = '  +0.43%  '

$ws.Range("E42").Value = '  +0.07%  '

$ws.Range("E43").Value = '  -0.37%  '

$ws.Range("D44").Value = '''13.50'
$ws.Range("E44").Value = '  +2.66%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '''6.05'
$ws.Range("E45").Value = '  +3.38%  '

$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").Value = '''0.0507'
$ws.Range("E46").Value = '  +1.63%  '

$ws.Range("D47").Value = '''1.07'
$ws.Range("E47").Value = '  -0.48%  '

$ws.Range("D48").Value = '0.0₆0135'
$ws.Range("E48").Value = '  -0.26%  '

$ws.Range("D49").Value = '''106.87'
$ws.Range("E49").Value = '  -2.09%  '

$ws.Range("D50").Value = '1.945.69'
$ws.Range("E50").Value = '  -0.20%  '

$ws.Range("E51").Value = '  +0.10%  '
